$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextCell "D2" "58.251.78"
$ws.Range("E2").Value = "  +0.48%  "

Set-TextCell "D3" "2.593.34"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.09%  "

Set-TextCell "D5" "522.79"
$ws.Range("E5").Value = "  +1.22%  "

Set-TextCell "D6" "144.54"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  +0.29%  "

Set-TextCell "D9" "2.615.73"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("E11").Value = "  -0.74%  "

Set-TextCell "D12" "0.334"
$ws.Range("E12").Value = "  -0.59%  "

$ws.Range("E13").Value = "  -0.70%  "

Set-TextCell "D14" "3.055.19"

Set-TextCell "D15" "58.203.74"
$ws.Range("E15").Value = "  +0.41%  "

Set-TextCell "D16" "20.56"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D17" "2.636.85"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D18" "0.0000134"
$ws.Range("E18").Value = "  -0.44%  "

Set-TextCell "D19" "340.45"
$ws.Range("E19").Value = "  +2.01%  "

$ws.Range("E20").Value = "  -0.35%  "

Set-TextCell "D21" "10.33"
$ws.Range("E21").Value = "  +0.22%  "

Set-TextCell "D22" "6.40"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").Value = "  -0.25%  "

Set-TextCell "D24" "65.54"
$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  -2.34%  "

Set-TextCell "D27" "2.715.80"
$ws.Range("E27").Value = "  -0.46%  "

Set-TextCell "D28" "0.996"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  -0.29%  "

Set-TextCell "D30" "0.0₃0752"
$ws.Range("E30").Value = "  -4.05%  "

$ws.Range("E32").Value = "  -5.09%  "

$ws.Range("E33").Value = "  +0.96%  "

Set-TextCell "D34" "18.88"
$ws.Range("E34").Value = "  +1.24%  "

Set-TextCell "D35" "149.86"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("E38").Value = "  -3.35%  "

Set-TextCell "D39" "0.848"
$ws.Range("E39").Value = "  +1.45%  "

Set-TextCell "D40" "1.46"
$ws.Range("E40").Value = "  +2.49%  "

Set-TextCell "D41" "36.11"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").Value = "  -0.95%  "

Set-TextCell "D43" "0.997"
$ws.Range("E43").Value = "  -0.27%  "

Set-TextCell "D44" "273.77"
$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("E47").Value = "  +0.41%  "

Set-TextCell "D48" "0.0525"
$ws.Range("E48").Value = "  -0.96%  "

Set-TextCell "D49" "18.82"
$ws.Range("E49").Value = "  -1.33%  "

Set-TextCell "D50" "19.16"
$ws.Range("E50").Value = "  +5.37%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D51" "1.981.61"
$ws.Range("E51").Value = "  -2.18%  "
